# Add a new "Bonus Question Solution" slide right after the existing slide,
# using the same "Title and Content" layout as slide 1 (layout index 2).
$p = $ppt.ActivePresentation
$s = $p.Slides.Add(2, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Bonus Question Solution"

# Body / content placeholder - build it paragraph by paragraph so each
# run keeps its own formatting (mirrors how the slide was authored).
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange
$tr.Text = "From 3x3 grid, we see the number"

$tr.InsertAfter("`n1, 4, 9") | Out-Null
$tr.InsertAfter("`nWhich is 1 x 1 + 2 x 2 + 3 x 3") | Out-Null
$tr.InsertAfter("`nVerify the pattern ") | Out-Null
$tr.InsertAfter("with 2x2 and 4x4 ") | Out-Null
$tr.InsertAfter("grid") | Out-Null
$tr.InsertAfter("`nSo, for 8x8 grid, we have") | Out-Null
$tr.InsertAfter("`n1 + 4 + 9 + 16 + 25 + 36 + 49 + 64 = 204") | Out-Null

# Indent the two "answer" lines one level in, now that every paragraph
# already exists (so the indent doesn't leak into later InsertAfter calls).
$tr.Paragraphs(2).IndentLevel = 2
$tr.Paragraphs(6).IndentLevel = 2
